$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wildo")

$ws.Range('A97').Value = 'Programador'
$ws.Range('A97').Font.Bold = $true
$ws.Range('B97').Value = 'Leois Linka'
$ws.Range('C97').Value = 'Tester'
$ws.Range('C97').Font.Bold = $true
$ws.Range('A98').Value = 'Modulo'
$ws.Range('A98').Font.Bold = $true
$ws.Range('B98').Value = 'Destinatarios'
$ws.Range('B98').Font.Bold = $true
$ws.Range('B99').Value = 'Tests'
$ws.Range('B99').Font.Bold = $true
$ws.Range('C99').Value = 'Resultado Esperado'
$ws.Range('C99').Font.Bold = $true
$ws.Range('D99').Value = 'Resultado Otenido'
$ws.Range("A100").Value = 1
$ws.Range('B100').Value = 'Ingresar un destinatario con con todos los campos obligatorios'
$ws.Range('C100').Value = 'Debe guardar el destinatario mostrarlo en la tabla inferior'
$ws.Range('C101').Value = 'y limpiar los campos del formulario'
$ws.Range("A103").Value = 2
$ws.Range('B103').Value = 'Registrar un cliente, sin ingresar todos los campos obligariorios'
$ws.Range('C103').Value = 'No debe de guardarse el cliente, los validadores deben '
$ws.Range('C104').Value = 'aparecer en rojo'
$ws.Range("A106").Value = 3
$ws.Range('B106').Value = 'Registrar un cliente sin ingresar al menos una direccion.'
$ws.Range('C106').Value = 'No debe guardar el destintario. Debe mostrar un mensaje informando que '
$ws.Range('C107').Value = 'al menos debe contener una direccion'
$ws.Range("A110").Value = 4
$ws.Range('B110').Value = 'Hacer click en agregar sin setear los campos correspondientes a la direccion'
$ws.Range('C110').Value = 'No debe agregar la direccion en la tabla de direcciones  a la derecha del '
$ws.Range('C111').Value = 'formulario'
$ws.Range("A113").Value = 5
$ws.Range('B113').Value = 'Eliminar un destinatario con referencias a otros modelos'
$ws.Range('C113').Value = 'No debe permitir eliminar un destinatario que esta siendo usado'
$ws.Range("A115").Value = 6
$ws.Range('B115').Value = 'Eliminar destinatarios sin referencias a otros modelos'
$ws.Range('C115').Value = 'Aparece un mensaje de confirmacion, preguntando si realmente desea'
$ws.Range('C116').Value = 'eliminar el destinatario'
$ws.Range("A118").Value = 7
$ws.Range('B118').Value = 'Ingresar muchos caracteres en los campos:'
$ws.Range('C118').Value = 'Nombre y Apellido: max cantidad de caracteres 50'
$ws.Range('C119').Value = 'Nro de Documento: max cantidad de caracteres 20'
$ws.Range('C120').Value = 'Ruc: max cantidad de caracteres 14'
$ws.Range('C121').Value = 'Lugar: max cantidad de caracteres 15'
$ws.Range('C122').Value = 'Ciudad: debe ser un autocomplete'
$ws.Range('C123').Value = 'Dirección:  texto'
$ws.Range('A125').Value = 'Modulo'
$ws.Range('A125').Font.Bold = $true
$ws.Range('B125').Value = 'Usuarios'
$ws.Range('B125').Font.Bold = $true
$ws.Range('B126').Value = 'Tests'
$ws.Range('B126').Font.Bold = $true
$ws.Range('C126').Value = 'Resultado Esperado'
$ws.Range('C126').Font.Bold = $true
$ws.Range('D126').Value = 'Resultado Otenido'
$ws.Range("A127").Value = 1
$ws.Range('B127').Value = 'Ingresar un usuario con con todos los campos obligatorios'
$ws.Range('C127').Value = 'Debe guardar el usuario, mostrarlo en la tabla inferior'
$ws.Range('C128').Value = 'y limpiar los campos del formulario'
$ws.Range("A130").Value = 2
$ws.Range('B130').Value = 'Registrar un usuario, sin ingresar todos los campos obligariorios'
$ws.Range('C130').Value = 'No debe de guardarse el cliente, las validaciones aparecen en rojo'
$ws.Range('C131').Value = 'No se habilita el boton guardar'
$ws.Range("A133").Value = 3
$ws.Range('B133').Value = 'Registrar un usuario sin ningun rol'
$ws.Range('C133').Value = 'No debe guardar el usuario. No se habilita el boton guardar'
$ws.Range("A140").Value = 4
$ws.Range('B140').Value = 'Eliminar un usuario con referencias a otros modelos'
$ws.Range('C140').Value = 'No debe permitir eliminar un usuario que esta siendo usado'
$ws.Range("A142").Value = 5
$ws.Range('B142').Value = 'Eliminar usuarios sin referencias a otros modelos'
$ws.Range('C142').Value = 'Aparece un mensaje de confirmacion, preguntando si realmente desea'
$ws.Range('C143').Value = 'eliminar el usuario'
$ws.Range("A145").Value = 6
$ws.Range('B145').Value = 'Ingresar muchos caracteres en los campos:'
$ws.Range('C145').Value = 'Empleado: autocomplete'
$ws.Range('C146').Value = 'Usuario: max cantidad de caracteres 20'
$ws.Range('C147').Value = 'Email: max cantidad de caracteres 50'
$ws.Range('C148').Value = 'Contraseña: no tiene limite'
$ws.Range('C149').Value = 'Confirmar contraseña: no tiene limite'
$ws.Range("A146").Select()
